# Commit: "updated order files with png"
#
# The stimulus order sheet referenced JPEG face/cue images; the experiment
# assets were switched to PNG, so every filename stored in the
# Filename_Left / Filename_Right columns (D and E) needs its extension
# updated from ".jpg" to ".png". Nothing else about those values changes.
#
# After updating the data, reset the view so the whole used range is
# selected/shown from the top, matching a fresh "select-all, scroll home"
# state rather than wherever the previous editor happened to leave the
# window scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    foreach ($col in @("D", "E")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val -like "*.jpg") {
            $cell.Value = ($val -replace "\.jpg$", ".png")
        }
    }
}

# Reset the saved view: select the full used range (no stray active-cell /
# scrolled-off top-left-cell left over from the previous session).
$usedRange.Select()

Write-Output "Renamed .jpg -> .png in Filename_Left/Filename_Right columns; reset selection."
